$d = $word.ActiveDocument

# The "Case" column is column 1 of the second table in the document.
$tbl = $d.Tables.Item(2)

# ---------------------------------------------------------------------
# 1) Update the case-label text in column 1 for each of the five rows.
#    Only the leading run of each cell needs editing, so setting the
#    cell Range's .Text replaces just that run's text in place while
#    leaving any trailing runs (e.g. spell-checked fragments) intact.
# ---------------------------------------------------------------------
$tbl.Cell(2,1).Range.Text = "1 – sea level rise of 2mm/year"
$tbl.Cell(3,1).Range.Text = "2 – As 1, with "
$tbl.Cell(4,1).Range.Text = "3 – historic changes + "
$tbl.Cell(5,1).Range.Text = "4 – As 3, with dredge in 2000*"
$tbl.Cell(6,1).Range.Text = "5 – As 4, with reclamation in 2020*"

# ---------------------------------------------------------------------
# 2) Give each of those five cells a full single-line border on all
#    four sides (0.25pt, automatic colour) instead of the plain shading
#    that was there before.
# ---------------------------------------------------------------------
$rowsToBorder = 2,3,4,5,6
foreach ($rowIdx in $rowsToBorder) {
    $cell = $tbl.Cell($rowIdx,1)
    $borders = $cell.Borders

    $borders.Item(-1).LineStyle = 1   # wdBorderTop    -> wdLineStyleSingle
    $borders.Item(-1).LineWidth = 2   # wdLineWidth025pt
    $borders.Item(-1).ColorIndex = 0  # wdAuto

    $borders.Item(-2).LineStyle = 1   # wdBorderLeft
    $borders.Item(-2).LineWidth = 2
    $borders.Item(-2).ColorIndex = 0

    $borders.Item(-3).LineStyle = 1   # wdBorderBottom
    $borders.Item(-3).LineWidth = 2
    $borders.Item(-3).ColorIndex = 0

    $borders.Item(-4).LineStyle = 1   # wdBorderRight
    $borders.Item(-4).LineWidth = 2
    $borders.Item(-4).ColorIndex = 0

    $borders.DistanceFromTop = 0
    $borders.DistanceFromLeft = 0
    $borders.DistanceFromBottom = 0
    $borders.DistanceFromRight = 0
}

Write-Output "edits applied"
